$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix implicit case issue: condition expression used lower-case "param1" which
#    didn't match the declared parameter name "Param1" -> fix casing.
$ws.Range("E5").Value = "contains(intArr, Param1)"

# 2. Add three new SmartRules test tables, built by copying the existing
#    "SmartMytestEasy" rule table (rows 25-28, columns D:F) which already has the
#    correct look & feel (borders, alignment, number format) and the right
#    "Cond1 / someAction / myRet" header row, then editing the title and pasting
#    it at the three new locations.

# Table 3: SmartMytestEasy3(Double PARAM1) - rows 39-42
$ws.Range("D25:F28").Copy($ws.Range("D39"))
$ws.Range("D39").Value = "SmartRules  Double[] SmartMytestEasy3(Double PARAM1)"

# Table 4: SmartMytestEasy4(Integer PARAM1) - rows 46-49
$ws.Range("D25:F28").Copy($ws.Range("D46"))
$ws.Range("D46").Value = "SmartRules  Double[] SmartMytestEasy4(Integer PARAM1)"

# Table 5: SmartMytestEasy5(Integer p1, Integer p2) - rows 54-57
$ws.Range("D25:F28").Copy($ws.Range("D54"))
$ws.Range("D54").Value = "SmartRules  Double[] SmartMytestEasy5(Integer p1, Integer p2)"

# Update the view to reflect where the user ended up after the edit.
$ws.Range("D3:G57")
$ws.Application.ActiveWindow.ScrollRow = 28
$ws.Range("J33").Select()
